# Update for new fish tagged and receiver download 4
# Adds a new deployment record (row 15) to Sheet1, reusing the existing
# date-format style from the row above it, and selects row 7 as the last
# active selection (matching the editing session in the source workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy A14's format (date number format, style index) down to A15 first,
# then overwrite its value with the new deployment date (3/1/2023).
$ws.Range("A14").Copy($ws.Range("A15"))
$ws.Range("A15").Value = 44986

$ws.Range("B15").Value = "SS"
$ws.Range("C15").Value = "MS"
$ws.Range("D15").Value = 138847
$ws.Range("E15").Value = "Mead Point #3"
$ws.Range("F15").Value = 27
$ws.Range("G15").Value = 31.254
$ws.Range("H15").Value = 27.5209
$ws.Range("I15").Value = 82
$ws.Range("J15").Value = 40.382
$ws.Range("K15").Value = 82.6730333

# Leave the same selection state recorded in the authored workbook
# (entire row 7 selected).
$ws.Rows(7).Select()
